$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists lines (line1..line6) followed by extraction points
# (extr1..extr8). This edit inserts two more lines ("line7", "line8")
# between the existing lines and the extraction points, which pushes the
# extraction-point rows down by two, and refreshes a handful of their
# C/D/E values (simulation results) along the way.
#
# Rows are written directly (rather than via EntireRow/Rows.Insert) so no
# incidental style entries get created; the bold/bordered "index column"
# look of column A is reproduced afterwards with a formats-only paste from
# an existing cell that already carries it.

# --- New row 8: line7 ---
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# --- New row 9: line8 ---
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# --- Row 10 (was extr1 on row 8): keep A/B, refresh C/D/E ---
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# --- Row 11 (was extr2 on row 9): keep A/B, refresh C/D/E ---
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# --- Row 12 (was extr3 on row 10): keep A/B, refresh C/D/E ---
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $false

# --- Row 13 (was extr4 on row 11): keep A/B, refresh C/D/E ---
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $true

# --- Row 14 (was extr5 on row 12): keep A/B, refresh C/D/E ---
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $true

# --- Row 15 (was extr6 on row 13): keep A/B, refresh C/D/E ---
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

# --- New row 16 (was extr7 on row 14): same C/D/E values, shifted down ---
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# --- New row 17 (was extr8 on row 15): same C/D/E values, shifted down ---
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true

# Re-apply the bold/bordered "index column" formatting (style of A2:A15) to
# the freshly written A8:A9 and A16:A17 cells.
$ws.Range("A2").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
